# The original template contained a field whose instrText read:
#   gd:for v| self.eAllStructuralFeatures
# The commit removes the dedicated "var" template construct ("gd:for v|")
# collapsing it down to a plain AQL expression; as typed live in Word this
# shows up as the literal text "gd" being replaced by "m" inside the field
# code, with Word's usual "_GoBack" last-edit bookmark landing right after
# the edited character. We reproduce that exact run layout.
#
# Because this runtime's Range/Find APIs do not address the hidden
# characters that make up a field's instrText (they always resolve to
# the very start of the story), we can't edit the field code in place with
# Find/Replace or Range.Text. Instead we locate the paragraph that owns
# the field, delete its contents, and re-insert the paragraph with the
# exact target run/bookmark layout via Range.InsertXML - the same net
# effect, expressed at the OOXML level the COM object model exposes
# through InsertXML.

$d = $word.ActiveDocument

$targetCode = "gd:for v| self.eAllStructuralFeatures "
$targetParagraph = $null

foreach ($f in $d.Fields) {
    if ($f.Code.Text -eq $targetCode) {
        $targetParagraph = $f.Code.Paragraphs.Item(1)
        break
    }
}

if ($targetParagraph -eq $null) {
    throw "Could not find the paragraph containing the 'gd:for v|...' field"
}

$paraRange = $targetParagraph.Range
$paraRange.Delete()

$insertionPoint = $targetParagraph.Range
$insertionPoint.Collapse(1)

$newParagraphXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:r><w:fldChar w:fldCharType="begin"/></w:r>' +
    '<w:r><w:instrText>m</w:instrText></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
    '<w:bookmarkEnd w:id="0"/>' +
    '<w:r><w:instrText xml:space="preserve">:for v| self.eAllStructuralFeatures </w:instrText></w:r>' +
    '<w:r><w:fldChar w:fldCharType="end"/></w:r>' +
    '<w:r><w:t>Un paragraphe tout nu</w:t></w:r>' +
    '</w:p>'

$insertionPoint.InsertXML($newParagraphXml)
